# Apply updated cryptocurrency price/volume values (Fri Jun 30 09:26:32 UTC 2023 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.800.61"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "'1.887.47"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'238.81"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.4766"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.2869"
$ws.Range("E8").Value = "  +4.86%  "
$ws.Range("D9").Value = "'0.06585"
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("D10").Value = "'18.73"
$ws.Range("E10").Value = "  +11.66%  "
$ws.Range("D11").Value = "'98.53"
$ws.Range("E11").Value = "  +16.98%  "
$ws.Range("D12").Value = "'1.890.43"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "'0.07600"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "'5.128"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").Value = "'0.6574"
$ws.Range("E15").Value = "  +5.54%  "
$ws.Range("D16").Value = "'307.20"
$ws.Range("E16").Value = "  +34.04%  "
$ws.Range("D17").Value = "'30.788.82"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "'13.17"
$ws.Range("E18").Value = "  +5.68%  "
$ws.Range("D19").Value = "'1.0000"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "'0.000007578"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("D21").Value = "'2.114.58"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'5.125"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").Value = "'6.170"
$ws.Range("E24").Value = "  +4.53%  "
$ws.Range("D25").Value = "'9.274"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").Value = "'167.24"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'20.30"
$ws.Range("E27").Value = "  +13.37%  "
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("E29").Value = "  +4.90%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").Value = "'4.176"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "'3.969"
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("D33").Value = "'0.05044"
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").Value = "'1.169"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").Value = "'0.7275"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").Value = "'2.715"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").Value = "'0.01946"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "'2.701"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "'2.072"
$ws.Range("E39").Value = "  +5.95%  "
$ws.Range("D40").Value = "'0.9022"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").Value = "'107.79"
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "'0.4203"
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("D44").Value = "'5.629"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").Value = "'7.350"
$ws.Range("E45").Value = "  +3.58%  "
$ws.Range("D46").Value = "'65.50"
$ws.Range("E46").Value = "  +6.56%  "
$ws.Range("D47").Value = "'9.059"
$ws.Range("E47").Value = "  +5.47%  "
$ws.Range("D48").Value = "'0.1224"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "'34.72"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").Value = "'0.05629"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").Value = "'1.391"
$ws.Range("E51").Value = "  +2.80%  "
